$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:K1) ---
$ws.Cells.Item(1,1).Value = "Row"
$ws.Cells.Item(1,2).Value = "Prognose"
$ws.Cells.Item(1,3).Value = "surveys"
$ws.Cells.Item(1,4).Value = "production"
$ws.Cells.Item(1,5).Value = "orders"
$ws.Cells.Item(1,6).Value = "turnover"
$ws.Cells.Item(1,7).Value = "financial"
$ws.Cells.Item(1,8).Value = "labor market"
$ws.Cells.Item(1,9).Value = "prices"
$ws.Cells.Item(1,10).Value = "national accounts"
$ws.Cells.Item(1,11).Value = "Revision"

# --- Column A date labels (A2:A12) ---
# Written as text: force Text number-format before the writes so Excel
# does not auto-convert the "yyyy-mm-dd" strings into date serials, then
# ClearFormats() to drop the now-unneeded custom style (keeps the cells on
# the default style, same as the source file).
$dateRange = $ws.Range("A2:A12")
$dateRange.NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "2025-03-30"
$ws.Cells.Item(3,1).Value = "2025-04-15"
$ws.Cells.Item(4,1).Value = "2025-04-30"
$ws.Cells.Item(5,1).Value = "2025-05-15"
$ws.Cells.Item(6,1).Value = "2025-05-30"
$ws.Cells.Item(7,1).Value = "2025-06-15"
$ws.Cells.Item(8,1).Value = "2025-06-30"
$ws.Cells.Item(9,1).Value = "2025-07-15"
$ws.Cells.Item(10,1).Value = "2025-07-30"
$ws.Cells.Item(11,1).Value = "2025-08-15"
$ws.Cells.Item(12,1).Value = "2025-08-30"
$dateRange.ClearFormats()

# --- Data cells (B2:K12) ---
$ws.Cells.Item(2,2).Value = 0.2906987142529091
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 0
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(2,10).Value = 0
$ws.Cells.Item(2,11).Value = 0

$ws.Cells.Item(3,2).Value = 0.34121197625062144
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0.011310990176631001
$ws.Cells.Item(3,5).Value = 0.0003464213005461224
$ws.Cells.Item(3,6).Value = 0.0030496675898262896
$ws.Cells.Item(3,7).Value = 0.010174587310235189
$ws.Cells.Item(3,8).Value = -0.00009230035772437486
$ws.Cells.Item(3,9).Value = 0.026057335300482827
$ws.Cells.Item(3,10).Value = 0
$ws.Cells.Item(3,11).Value = -0.0003334393222847454

$ws.Cells.Item(4,2).Value = 0.2933321662207311
$ws.Cells.Item(4,3).Value = -0.017104158159779077
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = -0.0017585313786715245
$ws.Cells.Item(4,6).Value = -0.0002934294269353389
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = -0.0020836370517838673
$ws.Cells.Item(4,9).Value = -0.027274543569465844
$ws.Cells.Item(4,10).Value = 0.0002461670252894271
$ws.Cells.Item(4,11).Value = 0.0003883225314558758

$ws.Cells.Item(5,2).Value = 0.28812702667753687
$ws.Cells.Item(5,3).Value = 0.057282445160034806
$ws.Cells.Item(5,4).Value = -0.023742352443871618
$ws.Cells.Item(5,5).Value = -0.00037726424627502943
$ws.Cells.Item(5,6).Value = 0.012384147202562912
$ws.Cells.Item(5,7).Value = -0.05717825491170156
$ws.Cells.Item(5,8).Value = -0.00017462230939389124
$ws.Cells.Item(5,9).Value = 0.006542518190742044
$ws.Cells.Item(5,10).Value = 0
$ws.Cells.Item(5,11).Value = 0.000058243814708092145

$ws.Cells.Item(6,2).Value = 0.42962805656564856
$ws.Cells.Item(6,3).Value = 0.19235491685294634
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 0.009202252069297605
$ws.Cells.Item(6,6).Value = 0.0006618688078815388
$ws.Cells.Item(6,7).Value = 0
$ws.Cells.Item(6,8).Value = 0.0006081582284800862
$ws.Cells.Item(6,9).Value = -0.0613590509545704
$ws.Cells.Item(6,10).Value = 0
$ws.Cells.Item(6,11).Value = 0.000032884884076489485

$ws.Cells.Item(7,2).Value = 0.3770294696911207
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = -0.06676569380682133
$ws.Cells.Item(7,5).Value = 0.0044190641709149735
$ws.Cells.Item(7,6).Value = -0.009567074871354279
$ws.Cells.Item(7,7).Value = 0.012600386232878935
$ws.Cells.Item(7,8).Value = 0
$ws.Cells.Item(7,9).Value = 0.003825970583703548
$ws.Cells.Item(7,10).Value = 0
$ws.Cells.Item(7,11).Value = 0.0028887608161503042

$ws.Cells.Item(8,2).Value = 0.2217036530838096
$ws.Cells.Item(8,3).Value = -0.1269016499594095
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = -0.013249514101887998
$ws.Cells.Item(8,6).Value = -0.013405984935775989
$ws.Cells.Item(8,7).Value = 0
$ws.Cells.Item(8,8).Value = 0.000027099098430844305
$ws.Cells.Item(8,9).Value = -0.0032516369111496284
$ws.Cells.Item(8,10).Value = 0
$ws.Cells.Item(8,11).Value = 0.0014558702024811687

$ws.Cells.Item(9,2).Value = -0.21689696575503065
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = -0.21710824991852384
$ws.Cells.Item(9,5).Value = 0.005715248225020182
$ws.Cells.Item(9,6).Value = -0.21587833737314274
$ws.Cells.Item(9,7).Value = -0.009309098831816166
$ws.Cells.Item(9,8).Value = -0.002613168975815504
$ws.Cells.Item(9,9).Value = -0.00023028997519626082
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,11).Value = 0.0008232780106341153

$ws.Cells.Item(10,2).Value = -0.051221557401531226
$ws.Cells.Item(10,3).Value = 0.1759154777981508
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 0.0006633848233777908
$ws.Cells.Item(10,6).Value = 0.00302357514904392
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = -0.0005863770230990097
$ws.Cells.Item(10,9).Value = 0.03859362307519536
$ws.Cells.Item(10,10).Value = -0.06367652851359518
$ws.Cells.Item(10,11).Value = 0.011742253044425749

$ws.Cells.Item(11,2).Value = 0.3889631566145159
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 0.4366445802498935
$ws.Cells.Item(11,5).Value = -0.011196862766811954
$ws.Cells.Item(11,6).Value = -0.03687622782930403
$ws.Cells.Item(11,7).Value = -0.031698183146630615
$ws.Cells.Item(11,8).Value = -0.0010641234879140346
$ws.Cells.Item(11,9).Value = 0.06281679858774117
$ws.Cells.Item(11,10).Value = 0
$ws.Cells.Item(11,11).Value = 0.021558732409073134

$ws.Cells.Item(12,2).Value = 0.3253066867255563
$ws.Cells.Item(12,3).Value = 0.02591966945706614
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = -0.000319629746019583
$ws.Cells.Item(12,6).Value = -0.0007118467749358623
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(12,8).Value = -0.0003810676207830869
$ws.Cells.Item(12,9).Value = -0.04220039686742383
$ws.Cells.Item(12,10).Value = 0
$ws.Cells.Item(12,11).Value = -0.04596319833686341

# --- Column width tweaks (best achievable on pixel-quantized grid) ---
$ws.Columns.Item(7).ColumnWidth = 14.3
$ws.Columns.Item(9).ColumnWidth = 15.3
